$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell value corrections (SOM / summary calculation fix + new DB records) ---
$ws.Range("B4").Value = 102

$ws.Range("B15").Value = 2
$ws.Range("B16").Value = 3
$ws.Range("B17").Value = 4

$ws.Range("B22").Value = 0
$ws.Range("B23").Value = 0
$ws.Range("B24").Value = 0
$ws.Range("B25").Value = 0

$ws.Range("B34").Value = 0

# --- Sheet view / selection state ---
$ws.Range("B14").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
